$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the top, pushing all existing data rows down by one.
$ws.Rows.Item(1).Insert()

# Populate the new header row with the CSV-style column names.
$ws.Cells.Item(1, 1).Value = "post_id"
$ws.Cells.Item(1, 2).Value = "post"
$ws.Cells.Item(1, 3).Value = "author"
$ws.Cells.Item(1, 4).Value = "topic"

# Turn on AutoFilter for the topic column (now D1:D201 after the insert).
$ws.Range("D1:D201").AutoFilter()

# Record the filter database as a hidden, sheet-scoped defined name (mirrors
# what Excel itself writes out when AutoFilter is applied).
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=postsArabic!D1")
$filterName.RefersTo = "=postsArabic!`$D`$1:`$D`$201"
$filterName.Visible = $false

# Restore the view: scroll so column D is visible and select E13, matching
# the author's final cursor position.
$ws.Range("D1").Select()
$ws.Range("E13").Select()
